$d = $word.ActiveDocument

# Helper: locate a unique substring in the document and return its [start,end)
# character offsets. Always re-searches from $d.Content so earlier edits can't
# leave a stale Find range/state behind.
function Find-Bounds([string]$text) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "text not found: $text" }
    return @($r.Start, $r.End)
}

# --- Step 1: add the Vue.js sentence as a new run inside the paragraph that is
# currently empty and immediately precedes "Split up the team creation...".
$splitBounds = Find-Bounds("Split up the team creation task")
$vueInsertPos = $splitBounds[0] - 1
$vueRange = $d.Range($vueInsertPos, $vueInsertPos)
$vueXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Discovered Vue.js as a potential alternative to JQM with custom elements. Seems to have a significant amount of compatibility with the project I will investigate this as an alternative.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$vueRange.InsertXML($vueXml) | Out-Null

# --- Step 2: replace " displaying all teams a user can see, actually making
# the team and a significant amount of generalisation to make the component
# reusable." (which currently spans two runs + the _GoBack bookmark) with the
# proofErr-wrapped three-run version. This also deletes the _GoBack bookmark
# from here; it is recreated at the end of the new content added in step 3.
$segBounds = Find-Bounds(" displaying all teams a user can see, actually making the team and a significant amount of generalisation to make the component reusable.")
$segRange = $d.Range($segBounds[0], $segBounds[1])
$splitXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> displaying all teams a user can see, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>actually making</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the team and a significant amount of generalisation to make the component reusable.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$segRange.InsertXML($splitXml) | Out-Null

# --- Step 3: append the fourteen new paragraphs (the 11/11-12/11 journal
# entries) right after the "Split up the team creation..." paragraph, ending
# with the _GoBack bookmark on the final ("Conversation with Jarod...") paragraph.
$tailBounds = Find-Bounds("the team and a significant amount of generalisation to make the component reusable.")
$tailRange = $d.Range($tailBounds[1], $tailBounds[1])
$newParasXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>12/11/19</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Developed the form mechanic further to allow for modification of document records using the same set of custom elements. Adding to the data properties and switch statements will allow any object to be modified by these elements. My placement </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>help</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> me significantly here, having programmatically generated a number of form elements using server side languages I understand some of the underlying properties and how to interact with them with JavaScript.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Found an issue with the basic query list generation class. It will work perfectly where you are accessing records directly in the query. The issue comes when you don’t want the query direction but a reference inside of it. The notification page implementation of this has an issue where it is applying listeners to the whole result. The solution to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>both of these</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> issues Is the same, I will make a new query list that applies listeners in a more manual way. This will give me control in implementing classes allowing me to apply listeners to references or not at all in the case of notifications.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Made changes to the how a user sees team data, because there will be more properties on a team in the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>future</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I think it no longer makes sense to have a separate local nickname on the user. This also </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>lead</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to UX issues where a user may think they are changing the whole group name when in fact it is only there local reference.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Conversation with Jarod increased my concern with the UX of the QR scanning process. This</w:t></w:r><w:r><w:t xml:space="preserve"> may become a reason for a user to not use the application and will be revisited in later sprints.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($newParasXml) | Out-Null

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
